# ---------------------------------------------------------------------------
# Applies the two content changes captured by the target diff:
#
#   1. The table on slide 5 switches its table style (tableStyleId) from
#      {E71A953F-41CB-4CF8-9629-E8E7F8D0F5EC} to
#      {E6D9FED5-4894-432A-8CFF-3CE58A470643}.
#
#   2. The deck's live theme (ppt/theme/theme2.xml -- the part actually
#      linked from the slide master / used by every slide) swaps its
#      "Integral / Red Violet" palette for the "Office Theme / Office"
#      palette that the (otherwise unused) theme1.xml part already carries.
#      We drive this the same way a user would from the Design tab /
#      Variants gallery: by editing each of the twelve theme colour slots
#      on the presentation's live colour scheme.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{E6D9FED5-4894-432A-8CFF-3CE58A470643}")
        }
    }
}

# --- 2. Theme colours -------------------------------------------------------
# Index : slot      : target (Office Theme) RGB  -> OLE (0x00BBGGRR) literal
#   1   : dk1        : 000000 -> 0x000000
#   2   : lt1        : FFFFFF -> 0xFFFFFF
#   3   : dk2        : 44546A -> 0x6A5444
#   4   : lt2        : E7E6E6 -> 0xE6E6E7
#   5   : accent1    : 5B9BD5 -> 0xD59B5B
#   6   : accent2    : ED7D31 -> 0x317DED
#   7   : accent3    : A5A5A5 -> 0xA5A5A5
#   8   : accent4    : FFC000 -> 0x00C0FF
#   9   : accent5    : 4472C4 -> 0xC47244
#  10   : accent6    : 70AD47 -> 0x47AD70
#  11   : hlink      : 0563C1 -> 0xC16305
#  12   : folHlink   : 954F72 -> 0x724F95

$targetColors = @(
    0x000000,
    0xFFFFFF,
    0x6A5444,
    0xE6E6E7,
    0xD59B5B,
    0x317DED,
    0xA5A5A5,
    0x00C0FF,
    0xC47244,
    0x47AD70,
    0xC16305,
    0x724F95
)

$colorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $targetColors[$i - 1]
}
